$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the two oldest dates (2025-09-21 and 2025-09-22) from the top of the
# data table -- deleting row 2 twice shifts every remaining row up by two,
# which reproduces the C-column re-alignment seen in the diff.
$ws.Rows.Item(2).Delete()
$ws.Rows.Item(2).Delete()

# Append the 7 newest days (2025-12-16 .. 2025-12-22) as fresh rows with
# zeroed Non-HTTPS/HTTPS counters, continuing straight after 2025-12-15
# (now row 85).
$newDates = @("2025-12-16", "2025-12-17", "2025-12-18", "2025-12-19", "2025-12-20", "2025-12-21", "2025-12-22")

$scratchRow = 1000
$row = 86
foreach ($d in $newDates) {
    # Writing a date-shaped string straight into .Value triggers Excel's
    # autodetect-as-date coercion. Build it as a text formula in a scratch
    # cell, then Copy/Paste it into place -- Copy preserves the literal
    # text (shared-string) type instead of reinterpreting it.
    $ws.Cells.Item($scratchRow, 1).Formula = '="' + $d + '"'
    $ws.Cells.Item($scratchRow, 1).Copy($ws.Cells.Item($row, 1))
    $ws.Cells.Item($scratchRow, 1).ClearContents()

    $ws.Cells.Item($row, 2).Value = 0.0
    $ws.Cells.Item($row, 3).Value = 0.0
    $row = $row + 1
}
